$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextCell "D2" '27.899.88'
Set-TextCell "E2" '  -2.12%  '
Set-TextCell "D3" '1.782.56'
Set-TextCell "E3" '  -2.00%  '
Set-TextCell "D4" '1.004'
Set-TextCell "E4" '  +0.13%  '
Set-TextCell "D5" '311.02'
Set-TextCell "E5" '  -1.41%  '
Set-TextCell "D6" '1.004'
Set-TextCell "E6" '  +0.15%  '
Set-TextCell "D7" '0.5116'
Set-TextCell "E7" '  +0.99%  '
Set-TextCell "D8" '0.3764'
Set-TextCell "E8" '  -1.88%  '
Set-TextCell "D9" '0.07764'
Set-TextCell "E9" '  -8.41%  '
Set-TextCell "D10" '41.33'
Set-TextCell "E10" '  -1.42%  '
Set-TextCell "D11" '1.083'
Set-TextCell "E11" '  -2.19%  '
Set-TextCell "D12" '1.004'
Set-TextCell "E12" '  +0.14%  '
Set-TextCell "D13" '6.177'
Set-TextCell "E13" '  -3.45%  '
Set-TextCell "D14" '20.16'
Set-TextCell "E14" '  -4.13%  '
Set-TextCell "D15" '1.774.39'
Set-TextCell "E15" '  -2.19%  '
Set-TextCell "D16" '7.178'
Set-TextCell "E16" '  -4.06%  '
Set-TextCell "D17" '91.91'
Set-TextCell "E17" '  -1.47%  '
Set-TextCell "D18" '0.00001066'
Set-TextCell "E18" '  -6.76%  '
Set-TextCell "D19" '0.06513'
Set-TextCell "E19" '  -2.82%  '
Set-TextCell "D20" '1.004'
Set-TextCell "E20" '  +0.13%  '
Set-TextCell "D21" '16.97'
Set-TextCell "E21" '  -4.09%  '
Set-TextCell "D22" '5.882'
Set-TextCell "E22" '  -3.08%  '
Set-TextCell "D23" '27.942.38'
Set-TextCell "E23" '  -2.05%  '
Set-TextCell "D25" '2.241'
Set-TextCell "E25" '  -1.26%  '
Set-TextCell "D26" '159.19'
Set-TextCell "E26" '  +0.21%  '
Set-TextCell "D27" '20.22'
Set-TextCell "E27" '  -4.79%  '
Set-TextCell "D28" '1.982.12'
Set-TextCell "E28" '  -1.95%  '
Set-TextCell "D29" '2.335'
Set-TextCell "E29" '  -1.97%  '
Set-TextCell "D30" '122.13'
Set-TextCell "E30" '  -3.06%  '
Set-TextCell "D31" '0.1069'
Set-TextCell "E31" '  -0.58%  '
Set-TextCell "D32" '1.035'
Set-TextCell "E32" '  -6.17%  '
Set-TextCell "D33" '3.639'
Set-TextCell "E33" '  -1.46%  '
Set-TextCell "D34" '5.469'
Set-TextCell "E34" '  -4.81%  '
Set-TextCell "D35" '0.07037'
Set-TextCell "E35" '  -4.31%  '
Set-TextCell "D36" '0.02302'
Set-TextCell "E36" '  -2.46%  '
Set-TextCell "D37" '0.2118'
Set-TextCell "E37" '  -4.74%  '
Set-TextCell "D38" '8.582'
Set-TextCell "E38" '  -1.50%  '
Set-TextCell "D41" '0.6078'
Set-TextCell "E41" '  -3.86%  '
Set-TextCell "D42" '1.149'
Set-TextCell "E42" '  -3.43%  '
Set-TextCell "D43" '1.332'
Set-TextCell "E43" '  -5.12%  '
Set-TextCell "D45" '13.01'
Set-TextCell "E45" '  -4.27%  '
Set-TextCell "D47" '126.60'
Set-TextCell "E47" '  +0.95%  '
Set-TextCell "D48" '1.215'
Set-TextCell "E48" '  +1.63%  '
Set-TextCell "D49" '1.886'
Set-TextCell "E49" '  -5.09%  '
Set-TextCell "D50" '0.06705'
Set-TextCell "E50" '  -4.13%  '
Set-TextCell "D51" '1.050'
Set-TextCell "E51" '  -1.75%  '
Set-TextCell "E24" '  -4.34%  '
Set-TextCell "B39" 'InternetComputer(DFINITY)'
Set-TextCell "C39" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell "D39" '4.999'
Set-TextCell "E39" '  -3.84%  '
Set-TextCell "B40" 'Aptos'
Set-TextCell "C40" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell "D40" '11.46'
Set-TextCell "E40" '  +1.88%  '
Set-TextCell "B44" 'Decentraland'
Set-TextCell "C44" 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell "D44" '0.5953'
Set-TextCell "E44" '  +0.90%  '
Set-TextCell "B46" 'PancakeSwap'
Set-TextCell "C46" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell "D46" '3.724'
Set-TextCell "E46" '  -0.68%  '
